$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 values: Day (date serial), 24 hourly prices (B:Y), Price_Daily_Avg (Z),
# Slot_4h_max (AA, unchanged text), Slot_4h_price (AB), Slot_2h_frist (AC, unchanged text),
# Slot_2h_frist_price (AD), Slot_2h_second (AE), Slot_2h_second_price (AF), Slot_min_price (AG)

$ws.Range("A2").Value = 45920

$ws.Range("B2").Value = 104.99
$ws.Range("C2").Value = 99.89
$ws.Range("D2").Value = 97
$ws.Range("E2").Value = 93.8
$ws.Range("F2").Value = 95.29000000000001
$ws.Range("G2").Value = 95.15000000000001
$ws.Range("H2").Value = 99.89
$ws.Range("I2").Value = 101
$ws.Range("J2").Value = 99.89
$ws.Range("K2").Value = 65.59999999999999
$ws.Range("L2").Value = 30
$ws.Range("M2").Value = 5.76
$ws.Range("N2").Value = 4.31
$ws.Range("O2").Value = 1.73
$ws.Range("P2").Value = 1.73
$ws.Range("Q2").Value = 0.65
$ws.Range("R2").Value = 4.31
$ws.Range("S2").Value = 15.43
$ws.Range("T2").Value = 55.08
$ws.Range("U2").Value = 93.56
$ws.Range("V2").Value = 110.36
$ws.Range("W2").Value = 110.91
$ws.Range("X2").Value = 104.6
$ws.Range("Y2").Value = 101

$ws.Range("Z2").Value = 66.33

$ws.Range("AB2").Value = 106.72
$ws.Range("AD2").Value = 110.64
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 102.8
$ws.Range("AG2").Value = "9h-18h"
